$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPDK Functions")
$ws.Activate()

# Remove the rows for tests that were dropped from the DPDK Functions sheet:
#   row 3 - PrepareFdbTableEntryforV4GeneveTunnel / dpdk_fdb_tx_geneve_test
#   row 6 - PrepareGeneveDecapModTableEntry / dpdk_geneve_decap_test
#   row 7 - PrepareGeneveEncapTableEntry / dpdk_geneve_encap_test
#   row 9 - PrepareVxlanDecapModTableEntry / dpdk_vxlan_decap_test
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(3).Delete()

# The remaining VxlanEncap row (now row 6) had its Base Class entry fixed
# from the "IpTunnelTest" outlier back to the standard "BaseTableTest".
$ws.Cells.Item(6, 5).Value = "BaseTableTest"

# Clear the leftover special (teal/italic/red) highlighting on the Test Name
# column and on the couple of other flagged cells so they match the plain
# formatting used by the rest of the table.
$ws.Cells.Item(3, 3).Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)

$ws.Cells.Item(2, 4).Copy()
$ws.Range("B2:C2").PasteSpecial(-4122)

$ws.Cells.Item(6, 3).Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("E6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the saved selection/cursor position recorded for this sheet.
$ws.Range("B9").Select() | Out-Null
